# 2.1.1.1e — add the 2023 data column (Q) to the table, mirroring the
# formatting already used for the 2022 column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source file had iterative calculation switched on; turn it back off.
$excel.Iteration = $false

# Bring column Q's formatting (number format, borders, fonts, ...) in line
# with column P before writing any values into it.
$ws.Range("P3:P14").Copy() | Out-Null
$ws.Range("Q3:Q14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Header (row 4) + the 2023 data points for each indicator (rows 5-14).
$ws.Range("Q4").Value = 2023
$ws.Range("Q5").Value = 74.605426356589135
$ws.Range("Q6").Value = 118.8
$ws.Range("Q7").Value = 71.61643835616438
$ws.Range("Q8").Value = 95.703125
$ws.Range("Q9").Value = 113.91018619934282
$ws.Range("Q10").Value = 108.21501014198785
$ws.Range("Q11").Value = 165.26684164479443
$ws.Range("Q12").Value = 48.504446240905416
$ws.Range("Q13").Value = 97.361348644026393
$ws.Range("Q14").Value = 52.747252747252752

# Row heights settled by Excel's auto-fit after the extra column was added.
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 16.5
$ws.Rows.Item(8).RowHeight = 16.5
$ws.Rows.Item(9).RowHeight = 16.5
$ws.Rows.Item(10).RowHeight = 16.5
$ws.Rows.Item(11).RowHeight = 16.5
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(13).RowHeight = 16.5
$ws.Rows.Item(14).RowHeight = 16.5

# Reset the saved selection back to the top-left cell.
$ws.Range("A1").Select() | Out-Null
